$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Make room for the new "2022-Q1" sheet right before "总计" by
#    temporarily renaming "总计" out of the way, copying the "2021-Q4"
#    sheet (same column layout as the new sheet) to sit right after
#    itself (i.e. right before the renamed "总计"), and renaming that
#    copy to "2022-Q1". Finally restore "总计"'s name.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Name = "总计_tmp"

$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$q4Sheet.Copy($null, $q4Sheet)

$ws2022 = $wb.Worksheets.Item("2021-Q4 (2)")
$ws2022.Name = "2022-Q1"

$wb.Worksheets.Item("总计_tmp").Name = "总计"

# ---------------------------------------------------------------------
# 2. The copied template (2021-Q4) has 16 rows (header + 15 funds); the
#    new 2022-Q1 data only has 13 funds, so drop the 2 extra rows.
# ---------------------------------------------------------------------
$ws2022.Rows("15:16").Delete()

# ---------------------------------------------------------------------
# 3. Overwrite the fund rows with the 2022-Q1 figures. Columns D-G hold
#    numeric-looking figures that are stored as text (matching the
#    source data), so they're entered with a leading apostrophe to force
#    text interpretation without altering the cell's number format.
# ---------------------------------------------------------------------
$ws2022.Range("A2").Value = 0
$ws2022.Range("B2").Value = "'009076"
$ws2022.Range("C2").Value = "工银瑞信圆兴混合"
$ws2022.Range("D2").Value = "'59.11"
$ws2022.Range("E2").Value = "'91.87"
$ws2022.Range("F2").Value = "'3.82"
$ws2022.Range("G2").Value = "'2.2580"
$ws2022.Range("H2").Value = 9
$ws2022.Range("A3").Value = 1
$ws2022.Range("B3").Value = "'010591"
$ws2022.Range("C3").Value = "富国中国中小盘混合(QDII)美元"
$ws2022.Range("D3").Value = "'35.75"
$ws2022.Range("E3").Value = "'86.53"
$ws2022.Range("F3").Value = "'5.52"
$ws2022.Range("G3").Value = "'1.9734"
$ws2022.Range("H3").Value = 1
$ws2022.Range("A4").Value = 2
$ws2022.Range("B4").Value = "'100061"
$ws2022.Range("C4").Value = "富国中国中小盘混合(QDII)人民币"
$ws2022.Range("D4").Value = "'35.75"
$ws2022.Range("E4").Value = "'86.53"
$ws2022.Range("F4").Value = "'5.52"
$ws2022.Range("G4").Value = "'1.9734"
$ws2022.Range("H4").Value = 1
$ws2022.Range("A5").Value = 3
$ws2022.Range("B5").Value = "'011006"
$ws2022.Range("C5").Value = "工银瑞信圆丰三年持有期混合"
$ws2022.Range("D5").Value = "'74.61"
$ws2022.Range("E5").Value = "'90.12"
$ws2022.Range("F5").Value = "'2.51"
$ws2022.Range("G5").Value = "'1.8727"
$ws2022.Range("H5").Value = 9
$ws2022.Range("A6").Value = 4
$ws2022.Range("B6").Value = "'005847"
$ws2022.Range("C6").Value = "富国沪港深业绩驱动混合A"
$ws2022.Range("D6").Value = "'44.03"
$ws2022.Range("E6").Value = "'74.38"
$ws2022.Range("F6").Value = "'2.71"
$ws2022.Range("G6").Value = "'1.1932"
$ws2022.Range("H6").Value = 8
$ws2022.Range("A7").Value = 5
$ws2022.Range("B7").Value = "'006752"
$ws2022.Range("C7").Value = "天弘港股通精选灵活配置混合A"
$ws2022.Range("D7").Value = "'6.86"
$ws2022.Range("E7").Value = "'85.37"
$ws2022.Range("F7").Value = "'8.25"
$ws2022.Range("G7").Value = "'0.5660"
$ws2022.Range("H7").Value = 1
$ws2022.Range("A8").Value = 6
$ws2022.Range("B8").Value = "'009029"
$ws2022.Range("C8").Value = "工银瑞信高质量成长混合A"
$ws2022.Range("D8").Value = "'15.47"
$ws2022.Range("E8").Value = "'88.00"
$ws2022.Range("F8").Value = "'2.95"
$ws2022.Range("G8").Value = "'0.4564"
$ws2022.Range("H8").Value = 9
$ws2022.Range("A9").Value = 7
$ws2022.Range("B9").Value = "'006753"
$ws2022.Range("C9").Value = "天弘港股通精选灵活配置混合C"
$ws2022.Range("D9").Value = "'2.49"
$ws2022.Range("E9").Value = "'85.37"
$ws2022.Range("F9").Value = "'8.25"
$ws2022.Range("G9").Value = "'0.2054"
$ws2022.Range("H9").Value = 1
$ws2022.Range("A10").Value = 8
$ws2022.Range("B10").Value = "'011117"
$ws2022.Range("C10").Value = "富国沪港深业绩驱动混合C"
$ws2022.Range("D10").Value = "'2.39"
$ws2022.Range("E10").Value = "'74.38"
$ws2022.Range("F10").Value = "'2.71"
$ws2022.Range("G10").Value = "'0.0648"
$ws2022.Range("H10").Value = 8
$ws2022.Range("A11").Value = 9
$ws2022.Range("B11").Value = "'009030"
$ws2022.Range("C11").Value = "工银瑞信高质量成长混合C"
$ws2022.Range("D11").Value = "'1.87"
$ws2022.Range("E11").Value = "'88.00"
$ws2022.Range("F11").Value = "'2.95"
$ws2022.Range("G11").Value = "'0.0552"
$ws2022.Range("H11").Value = 9
$ws2022.Range("A12").Value = 10
$ws2022.Range("B12").Value = "'011969"
$ws2022.Range("C12").Value = "建信港股通精选混合A"
$ws2022.Range("D12").Value = "'1.01"
$ws2022.Range("E12").Value = "'57.64"
$ws2022.Range("F12").Value = "'4.01"
$ws2022.Range("G12").Value = "'0.0405"
$ws2022.Range("H12").Value = 8
$ws2022.Range("A13").Value = 11
$ws2022.Range("B13").Value = "'005259"
$ws2022.Range("C13").Value = "建信龙头企业股票"
$ws2022.Range("D13").Value = "'1.04"
$ws2022.Range("E13").Value = "'83.45"
$ws2022.Range("F13").Value = "'3.06"
$ws2022.Range("G13").Value = "'0.0318"
$ws2022.Range("H13").Value = 8
$ws2022.Range("A14").Value = 12
$ws2022.Range("B14").Value = "'011970"
$ws2022.Range("C14").Value = "建信港股通精选混合C"
$ws2022.Range("D14").Value = "'0.33"
$ws2022.Range("E14").Value = "'57.64"
$ws2022.Range("F14").Value = "'4.01"
$ws2022.Range("G14").Value = "'0.0132"
$ws2022.Range("H14").Value = 8

# ---------------------------------------------------------------------
# 4. Prepend a "2022-Q1" summary row to the "总计" sheet: insert a new
#    row 2 (pushing the old rows down), copy the index-column (A) style
#    from the row that used to be in that slot, fill in the new figures,
#    and re-number the old rows' index column (which counts up from 0).
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Rows("2:2").Insert()
$wsTotal.Range("A2:D2").ClearFormats()

$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 13
$wsTotal.Range("D2").Value = 10.7

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("A5").Value = 3
$wsTotal.Range("A6").Value = 4
$wsTotal.Range("A7").Value = 5
